$wb = $excel.ActiveWorkbook

# Sheet 1: "Means"
$wsMeans = $wb.Worksheets.Item("Means")

# Rename header labels (shared strings used by both sheets)
$wsMeans.Range("B1").Value = "Rural Areas (National Average)"
$wsMeans.Range("C1").Value = "Rural Areas (State Average)"

# Updated mean values in column B (rows 2-10)
$wsMeans.Range("B2").Value = 84
$wsMeans.Range("B3").Value = 7.6
$wsMeans.Range("B4").Value = 8.2
$wsMeans.Range("B5").Value = 10
$wsMeans.Range("B6").Value = 67
$wsMeans.Range("B7").Value = 6.8
$wsMeans.Range("B8").Value = 5.1
$wsMeans.Range("B9").Value = 26
$wsMeans.Range("B10").Value = 0.32

# Sheet 2: "Standard Deviations"
$wsSD = $wb.Worksheets.Item("Standard Deviations")

# Rename header labels (shared strings used by both sheets)
$wsSD.Range("B1").Value = "Rural Areas (National Average) SD"
$wsSD.Range("C1").Value = "Rural Areas (State Average) SD"

# Updated standard deviation values in column B (rows 2-9)
$wsSD.Range("B2").Value = 19
$wsSD.Range("B3").Value = 16
$wsSD.Range("B4").Value = 12
$wsSD.Range("B5").Value = 15
$wsSD.Range("B6").Value = 28
$wsSD.Range("B7").Value = 7.6
$wsSD.Range("B8").Value = 6.4
$wsSD.Range("B9").Value = 8.6
